# Apply updated dSF (column F) values for the dunning_dane sheet.
# These represent a repull of data / push of all data with a mean
# calculation update, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 1
    3  = 2
    4  = -2
    5  = -3
    6  = -2
    7  = 2
    8  = 5
    9  = 1
    10 = 4
    11 = 0
    12 = -2
    14 = 2
    15 = -1
    16 = -1
    17 = 0
    18 = 3
    19 = -4
    21 = 3
    22 = -5
    23 = -1
    24 = -3
    25 = -1
    26 = 2
    28 = -1
    29 = 2
    31 = -2
    33 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
